$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.323.81"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").Value = "2.980.40"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Formula = '="565.53"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Formula = '="137.65"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = "  +3.24%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "2.975.15"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Formula = '="0.132"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").Formula = '="5.37"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E11").Value = "  +11.27%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Formula = '="33.63"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "3.476.16"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Formula = '="7.06"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "2.983.00"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "59.358.57"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Formula = '="436.46"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E20").Value = "  +4.75%  "
$ws.Range("D21").Formula = '="13.57"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").Formula = '="0.720"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = "  +2.63%  "
$ws.Range("D23").Formula = '="7.01"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Formula = '="13.19"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Formula = '="79.81"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D27").Formula = '="2.22"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = "  +9.17%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").Formula = '="7.71"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Formula = '="6.23"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("D32").Formula = '="25.73"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  +7.81%  "
$ws.Range("D34").Formula = '="0.0"&UNICHAR(8323)&"0767"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E34").Value = "  +9.57%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Formula = '="0.986"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Formula = '="5.88"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Formula = '="48.62"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Formula = '="8.69"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Formula = '="2.75"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("D41").Formula = '="399.52"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "2.730.49"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").Formula = '="0.250"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("D47").Formula = '="34.50"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E47").Value = "  +18.13%  "
$ws.Range("D48").Formula = '="122.18"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").Formula = '="1.99"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Formula = '="23.15"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E51").Value = "  +1.49%  "
